# Update cryptos list with latest price/volume figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.892.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.29%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.483.72"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.15"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.630"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.32%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.635"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "

$ws.Range("E10").Value = "  +4.31%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "55.84"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000277"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.72%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.23"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.25%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.044.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.498.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.23"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.903.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.64%  "

$ws.Range("E20").Value = "  +1.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "409.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.16%  "

$ws.Range("E22").Value = "  +7.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +8.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.37%  "

$ws.Range("E27").Value = "  -2.20%  "

$ws.Range("E28").Value = "  -2.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.28%  "

$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.70"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.68%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "586.05"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -9.47%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.92%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.90"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.88%  "

$ws.Range("E36").Value = "  +1.52%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0791"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.383"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.228.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.55%  "

$ws.Range("E43").Value = "  +0.27%  "

$ws.Range("E44").Value = "  +2.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.52%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.54"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.22%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0418"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.62%  "

$ws.Range("E48").Value = "  -5.98%  "

$ws.Range("E49").Value = "  +0.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.01"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.52%  "
